$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 - this shifts all existing rows (6..18) down to (7..19)
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the latest weekly entry
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44614
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100103
$ws.Range("H6").Value = "Frutos de hueso (carozo)"
$ws.Range("I6").Value = 100103002
$ws.Range("J6").Value = "Ciruela"
$ws.Range("K6").Value = "Angeleno"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 19000
$ws.Range("P6").Value = 18500
$ws.Range("Q6").Value = "$/bandeja 18 kilos granel"
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 1028
$ws.Range("T6").Value = 18
